$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 90, shifting existing rows 90..199 down to 91..200
$ws.Rows.Item(90).Insert()

# Fill in the new record for row 90
$ws.Cells.Item(90, 1).Value = 9
$ws.Cells.Item(90, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(90, 3).Value = "Metropolitana"
$ws.Cells.Item(90, 4).Value = 44482
$ws.Cells.Item(90, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(90, 5).Value = 13
$ws.Cells.Item(90, 6).Value = 100112052
$ws.Cells.Item(90, 7).Value = "Albahaca"
$ws.Cells.Item(90, 8).Value = "Sin especificar"
$ws.Cells.Item(90, 9).Value = "Primera"
$ws.Cells.Item(90, 10).Value = 43
$ws.Cells.Item(90, 11).Value = 6000
$ws.Cells.Item(90, 12).Value = 7000
$ws.Cells.Item(90, 13).Value = 6512
$ws.Cells.Item(90, 14).Value = "`$/docena de matas"
$ws.Cells.Item(90, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(90, 16).Value = 1085
$ws.Cells.Item(90, 17).Value = 6
$ws.Cells.Item(90, 18).Value = "Hortaliza"
